$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# New time-registration entry in row 19
$ws.Range("A19").Value = "Lav SD0101 og DCD0101"
$ws.Range("B19").Value = "Software Architect"
$ws.Range("C19").Value2 = 43888
$ws.Range("D19").Value2 = 0.52083333333333337
$ws.Range("E19").Value2 = 0.66666666666666663

# Update the active selection to the new row's date cell
$ws.Range("C19").Select()
